$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 15003331
$ws.Range("J40").Value = 42860064
$ws.Range("L40").Value = 42860064
$ws.Range("N40").Value = -42860414
$ws.Range("H62").Value = 76927624
$ws.Range("J62").Value = 83339520
$ws.Range("L62").Value = 83339520
$ws.Range("N62").Value = -83340768
$ws.Range("H64").Value = 3088.5
$ws.Range("J64").Value = 4003
$ws.Range("L64").Value = 4003
$ws.Range("N64").Value = -4499
$ws.Range("H65").Value = 76927624
$ws.Range("J65").Value = 83339520
$ws.Range("L65").Value = 416697600
$ws.Range("N65").Value = -416703840
$ws.Range("H67").Value = 3088.5
$ws.Range("J67").Value = 4003
$ws.Range("L67").Value = 4003
$ws.Range("N67").Value = -5719
$ws.Range("H74").Value = 3195
$ws.Range("J74").Value = 3195
$ws.Range("L74").Value = 3195
$ws.Range("N74").Value = -5067
$ws.Range("H77").Value = 3195
$ws.Range("J77").Value = 3195
$ws.Range("L77").Value = 15975
$ws.Range("N77").Value = -25335
$ws.Range("H98").Value = 3606415.5
$ws.Range("I98").Value = 4547749.5
$ws.Range("K98").Value = 4547749.5
$ws.Range("M98").Value = -4546251.5
$ws.Range("H122").Value = 3606415.5
$ws.Range("I122").Value = 4547749.5
$ws.Range("K122").Value = 13643248.5
$ws.Range("M122").Value = -13640798.5
$ws.Range("H125").Value = 1399.25
$ws.Range("H132").Value = 3687.879
$ws.Range("I132").Value = 3339.0715
$ws.Range("J132").Value = 5641.2
$ws.Range("K132").Value = 10017.2145
$ws.Range("L132").Value = 16923.6
$ws.Range("M132").Value = -7487.2145
$ws.Range("N132").Value = -21983.6
$ws.Range("H135").Value = 692.6389
$ws.Range("I135").Value = 654.625
$ws.Range("K135").Value = 5891.625
$ws.Range("M135").Value = -3356.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1552.826
$ws.Range("I2").Value = 1464.4706
$ws.Range("J2").Value = 1803.1666
$ws.Range("K2").Value = 1464.4706
$ws.Range("L2").Value = 1803.1666
$ws.Range("M2").Value = -1351.4706
$ws.Range("N2").Value = -2029.1666
$ws.Range("H19").Value = 3749.5
$ws.Range("I19").Value = 3749.5
$ws.Range("K19").Value = 3749.5
$ws.Range("M19").Value = -3520.5
$ws.Range("H61").Value = 4993.1
$ws.Range("I61").Value = 3066.375
$ws.Range("K61").Value = 3066.375
$ws.Range("M61").Value = -2854.375
$ws.Range("H116").Value = 1552.826
$ws.Range("I116").Value = 1464.4706
$ws.Range("J116").Value = 1803.1666
$ws.Range("K116").Value = 1464.4706
$ws.Range("L116").Value = 1803.1666
$ws.Range("M116").Value = 829.5293999999999
$ws.Range("N116").Value = -6391.1666
$ws.Range("H132").Value = 1943.5306
$ws.Range("I132").Value = 1937.5312
$ws.Range("J132").Value = 1954.8235
$ws.Range("K132").Value = 5812.5936
$ws.Range("L132").Value = 5864.470499999999
$ws.Range("M132").Value = -3282.5936
$ws.Range("N132").Value = -10924.4705
$ws.Range("H136").Value = 4993.1
$ws.Range("I136").Value = 3066.375
$ws.Range("K136").Value = 9199.125
$ws.Range("M136").Value = -6649.125
$ws.Range("H140").Value = 61061.285
$ws.Range("J140").Value = 61061.285
$ws.Range("L140").Value = 61061.285
$ws.Range("N140").Value = -71421.285
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1552.826
$ws.Range("I3").Value = 1464.4706
$ws.Range("J3").Value = 1803.1666
$ws.Range("K3").Value = 1464.4706
$ws.Range("L3").Value = 1803.1666
$ws.Range("M3").Value = -1350.4706
$ws.Range("N3").Value = -2031.1666
$ws.Range("H20").Value = 1932
$ws.Range("I20").Value = 1718
$ws.Range("J20").Value = 2176.5715
$ws.Range("K20").Value = 1718
$ws.Range("L20").Value = 2176.5715
$ws.Range("M20").Value = -1471
$ws.Range("N20").Value = -2670.5715
$ws.Range("H25").Value = 2062.5
$ws.Range("I25").Value = 2062.5
$ws.Range("K25").Value = 2062.5
$ws.Range("M25").Value = -1827.5
$ws.Range("H105").Value = 2852.889
$ws.Range("I105").Value = 2696.7144
$ws.Range("K105").Value = 2696.7144
$ws.Range("M105").Value = -949.7143999999998
$ws.Range("H134").Value = 2006.238
$ws.Range("I134").Value = 1618.5555
$ws.Range("K134").Value = 4855.666499999999
$ws.Range("M134").Value = -2320.666499999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 272120.47
$ws.Range("I31").Value = 371672.56
$ws.Range("K31").Value = 371672.56
$ws.Range("M31").Value = -371377.56
$ws.Range("H34").Value = 272120.47
$ws.Range("I34").Value = 371672.56
$ws.Range("K34").Value = 371672.56
$ws.Range("M34").Value = -371470.56
$ws.Range("H132").Value = 2732.4119
$ws.Range("I132").Value = 2644.1304
$ws.Range("K132").Value = 7932.3912
$ws.Range("M132").Value = -5402.3912
$ws.Range("H134").Value = 9989.733
$ws.Range("I134").Value = 13309.1
$ws.Range("K134").Value = 39927.3
$ws.Range("M134").Value = -37392.3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15386073
$ws.Range("J131").Value = 1573.8889
$ws.Range("L131").Value = 4721.6667
$ws.Range("N131").Value = -14801.6667
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 50330
$ws.Range("I6").Value = 50330
$ws.Range("K6").Value = 50330
$ws.Range("M6").Value = -50217
$ws.Range("H16").Value = 50330
$ws.Range("I16").Value = 50330
$ws.Range("K16").Value = 50330
$ws.Range("M16").Value = -50080
$ws.Range("H107").Value = 2436.1
$ws.Range("I107").Value = 2495.8572
$ws.Range("J107").Value = 2296.6667
$ws.Range("K107").Value = 2495.8572
$ws.Range("L107").Value = 2296.6667
$ws.Range("M107").Value = -575.8571999999999
$ws.Range("N107").Value = -6136.6667
$ws.Range("H126").Value = 6527.5293
$ws.Range("I126").Value = 8569.299999999999
$ws.Range("J126").Value = 3610.7144
$ws.Range("K126").Value = 25707.9
$ws.Range("L126").Value = 10832.1432
$ws.Range("M126").Value = -23237.9
$ws.Range("N126").Value = -15772.1432
$ws.Range("H132").Value = 24826.041
$ws.Range("I132").Value = 26667.092
$ws.Range("K132").Value = 80001.276
$ws.Range("M132").Value = -77471.276
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7426.864
$ws.Range("I7").Value = 7414.8
$ws.Range("K7").Value = 7414.8
$ws.Range("M7").Value = -7302.8
$ws.Range("H20").Value = 7332.5
$ws.Range("I20").Value = 4665
$ws.Range("K20").Value = 4665
$ws.Range("M20").Value = -4439
$ws.Range("H22").Value = 828.2143
$ws.Range("I22").Value = 623.5714
$ws.Range("K22").Value = 623.5714
$ws.Range("M22").Value = -328.5714
$ws.Range("H27").Value = 828.2143
$ws.Range("I27").Value = 623.5714
$ws.Range("K27").Value = 623.5714
$ws.Range("M27").Value = -516.5714
$ws.Range("H40").Value = 5456
$ws.Range("I40").Value = 5570
$ws.Range("K40").Value = 5570
$ws.Range("M40").Value = -5434
$ws.Range("H126").Value = 7426.864
$ws.Range("I126").Value = 7414.8
$ws.Range("K126").Value = 22244.4
$ws.Range("M126").Value = -19774.4
$ws.Range("H136").Value = 1220.8125
$ws.Range("I136").Value = 1098.76
$ws.Range("K136").Value = 3296.28
$ws.Range("M136").Value = -746.2799999999997
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 20428.285
$ws.Range("I14").Value = 19399.6
$ws.Range("J14").Value = 23000
$ws.Range("K14").Value = 19399.6
$ws.Range("L14").Value = 23000
$ws.Range("M14").Value = -19231.6
$ws.Range("N14").Value = -23336
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H51").Value = 11824.523
$ws.Range("I51").Value = 13394.467
$ws.Range("K51").Value = 13394.467
$ws.Range("M51").Value = -12884.467
$ws.Range("H62").Value = 9872
$ws.Range("I62").Value = 5999.3335
$ws.Range("J62").Value = 12195.6
$ws.Range("K62").Value = 5999.3335
$ws.Range("L62").Value = 12195.6
$ws.Range("M62").Value = -5375.3335
$ws.Range("N62").Value = -13443.6
$ws.Range("H65").Value = 9872
$ws.Range("I65").Value = 5999.3335
$ws.Range("J65").Value = 12195.6
$ws.Range("K65").Value = 29996.6675
$ws.Range("L65").Value = 60978
$ws.Range("M65").Value = -26876.6675
$ws.Range("N65").Value = -67218
$ws.Range("H80").Value = 29998
$ws.Range("I80").Value = 29998
$ws.Range("K80").Value = 29998
$ws.Range("M80").Value = -29000
$ws.Range("H83").Value = 29998
$ws.Range("I83").Value = 29998
$ws.Range("K83").Value = 89994
$ws.Range("M83").Value = -85002
$ws.Range("H107").Value = 986.7692
$ws.Range("I107").Value = 863.625
$ws.Range("K107").Value = 2590.875
$ws.Range("M107").Value = -670.875
$ws.Range("H126").Value = 3111
$ws.Range("I126").Value = 2905.5
$ws.Range("K126").Value = 8716.5
$ws.Range("M126").Value = -6246.5
$ws.Range("H132").Value = 1603.7273
$ws.Range("I132").Value = 1564.1
$ws.Range("K132").Value = 4692.299999999999
$ws.Range("M132").Value = -2162.299999999999
